$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.119.89'
Set-TextValue $ws.Range("E2") '  +1.17%  '
Set-TextValue $ws.Range("D3") '1.768.39'
Set-TextValue $ws.Range("E3") '  +1.32%  '
Set-TextValue $ws.Range("E4") '  +0.07%  '
Set-TextValue $ws.Range("D5") '238.08'
Set-TextValue $ws.Range("E5") '  +0.29%  '
Set-TextValue $ws.Range("E6") '  +0.09%  '
Set-TextValue $ws.Range("D7") '0.5239'
Set-TextValue $ws.Range("E7") '  +3.69%  '
Set-TextValue $ws.Range("D8") '0.2765'
Set-TextValue $ws.Range("E8") '  +4.37%  '
Set-TextValue $ws.Range("D9") '40.60'
Set-TextValue $ws.Range("E9") '  -3.19%  '
Set-TextValue $ws.Range("D10") '0.06222'
Set-TextValue $ws.Range("E10") '  +1.01%  '
Set-TextValue $ws.Range("D11") '1.774.75'
Set-TextValue $ws.Range("E11") '  +1.80%  '
Set-TextValue $ws.Range("D12") '15.97'
Set-TextValue $ws.Range("E12") '  +3.59%  '
Set-TextValue $ws.Range("D13") '0.07024'
Set-TextValue $ws.Range("E13") '  +1.51%  '
Set-TextValue $ws.Range("D14") '0.6499'
Set-TextValue $ws.Range("E14") '  +8.44%  '
Set-TextValue $ws.Range("D15") '4.528'
Set-TextValue $ws.Range("E15") '  +0.70%  '
Set-TextValue $ws.Range("D16") '78.43'
Set-TextValue $ws.Range("E16") '  +1.96%  '
Set-TextValue $ws.Range("E17") '  +0.08%  '
Set-TextValue $ws.Range("D18") '1.0000'
Set-TextValue $ws.Range("E18") '  +0.05%  '
Set-TextValue $ws.Range("D19") '26.125.74'
Set-TextValue $ws.Range("E19") '  +1.15%  '
Set-TextValue $ws.Range("D20") '11.76'
Set-TextValue $ws.Range("E20") '  +1.34%  '
Set-TextValue $ws.Range("D21") '0.000006789'
Set-TextValue $ws.Range("E21") '  -0.90%  '
Set-TextValue $ws.Range("D22") '2.000.75'
Set-TextValue $ws.Range("E22") '  +1.83%  '
Set-TextValue $ws.Range("D23") '4.090'
Set-TextValue $ws.Range("E23") '  +1.27%  '
Set-TextValue $ws.Range("D24") '8.439'
Set-TextValue $ws.Range("E24") '  +3.76%  '
Set-TextValue $ws.Range("D25") '5.211'
Set-TextValue $ws.Range("E25") '  +0.17%  '
Set-TextValue $ws.Range("D26") '138.04'
Set-TextValue $ws.Range("E26") '  +0.09%  '
Set-TextValue $ws.Range("D27") '1.489'
Set-TextValue $ws.Range("E27") '  -1.68%  '
Set-TextValue $ws.Range("D28") '1.869'
Set-TextValue $ws.Range("E28") '  +3.44%  '
Set-TextValue $ws.Range("D29") '15.20'
Set-TextValue $ws.Range("E29") '  +1.40%  '
Set-TextValue $ws.Range("D30") '102.78'
Set-TextValue $ws.Range("E30") '  -0.65%  '
Set-TextValue $ws.Range("D31") '0.08405'
Set-TextValue $ws.Range("E31") '  +3.71%  '
Set-TextValue $ws.Range("D32") '3.737'
Set-TextValue $ws.Range("E32") '  -0.80%  '
Set-TextValue $ws.Range("D33") '3.472'
Set-TextValue $ws.Range("E33") '  +0.07%  '
Set-TextValue $ws.Range("D34") '0.04460'
Set-TextValue $ws.Range("E34") '  -1.26%  '
Set-TextValue $ws.Range("D35") '2.656'
Set-TextValue $ws.Range("E35") '  +0.17%  '
Set-TextValue $ws.Range("D36") '1.009'
Set-TextValue $ws.Range("E36") '  +2.59%  '
Set-TextValue $ws.Range("D37") '0.6115'
Set-TextValue $ws.Range("E37") '  +0.71%  '
Set-TextValue $ws.Range("D38") '2.768'
Set-TextValue $ws.Range("E38") '  +3.60%  '
Set-TextValue $ws.Range("D39") '1.995'
Set-TextValue $ws.Range("E39") '  +4.43%  '
Set-TextValue $ws.Range("D40") '0.01587'
Set-TextValue $ws.Range("E40") '  +2.58%  '
Set-TextValue $ws.Range("E41") '  +0.30%  '
Set-TextValue $ws.Range("D42") '103.06'
Set-TextValue $ws.Range("E42") '  +0.00%  '
Set-TextValue $ws.Range("D43") '0.3885'
Set-TextValue $ws.Range("E43") '  +2.01%  '
Set-TextValue $ws.Range("D44") '0.7525'
Set-TextValue $ws.Range("E44") '  +2.89%  '
Set-TextValue $ws.Range("D45") '4.954'
Set-TextValue $ws.Range("E45") '  -2.66%  '
Set-TextValue $ws.Range("B46") 'Cronos'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D46") '0.05514'
Set-TextValue $ws.Range("E46") '  +3.09%  '
Set-TextValue $ws.Range("B47") 'Aptos'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D47") '6.451'
Set-TextValue $ws.Range("E47") '  +9.06%  '
Set-TextValue $ws.Range("E48") '  +0.85%  '
Set-TextValue $ws.Range("D49") '30.38'
Set-TextValue $ws.Range("E49") '  +0.59%  '
Set-TextValue $ws.Range("D50") '52.99'
Set-TextValue $ws.Range("E50") '  +0.78%  '
Set-TextValue $ws.Range("D51") '0.3476'
Set-TextValue $ws.Range("E51") '  +0.47%  '
